$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F17").Value = -8
$ws.Range("F18").Value = 5
$ws.Range("F22").Value = -7
$ws.Range("F23").Value = -7
$ws.Range("F24").Value = -8
